# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output figures.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 282
    4  = 10827
    5  = 9742
    6  = 584
    8  = 701
    9  = 97
    10 = 9
    13 = 9518
    15 = 2427
    16 = 33
    17 = 72
    18 = 374
    19 = 10833
    20 = 10746
    22 = 12
    25 = 13
    26 = 10
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
